# Ending 2018 exclusion identification after roll correction.
# Adds a new batch of "Exclusions after roll correction" entries (column H)
# for rows 144-194, marks one already-entered label (E155) with a distinct
# font, and updates the active view/selection to where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H entries: "Exclusions after roll correction" -------------
# Column H already holds this kind of data for earlier rows (1-141); we are
# simply continuing to fill it in down to the end of the log (row 194).
# These values are short numeric-looking range strings (e.g. "1720-",
# "-240 1085-1315 ...") that must stay text, so we force a text number
# format before assigning them.

$hValues = [ordered]@{
    144 = "-7000"
    147 = "5680-"
    148 = "4330-"
    149 = "130-310 830-1570 1760-1860 3580-3680"
    152 = "-570 450-580 730-840 970-1170 1370-1510 1630-1820 1960-2080 2260-2400 2545-2690 3190-3300 3560-3660 3850-3980 4150-4320 4460-4650"
    153 = "920-1240 1580-1880 2820-4150 4490-4580 4750-4920 5020-5240 5340-5540 5640-6160 6610-10810 10910-11030 11280-11430"
    155 = "1180-1760 4830-5000"
    156 = "8200-"
    158 = "1540-1630 2120-2170"
    159 = "855-1350 1600-1720 2060-2210 2320-2440 2720-2880 3830-3970 10180-10290 10470-10530 13100-13160 15250-15300 15380-15430"
    160 = "1300-1410 1600-1740"
    162 = "7850-8000 8170-8310 9400-9600 9780-9870 19780-19950 21420-21600"
    163 = "1320-1400"
    164 = "-930"
    170 = "-355 1100-1200"
    172 = "100-150"
    173 = "5430-"
    174 = "1460-1580"
    177 = "6175-6275"
    178 = "-320"
    180 = "-100 630-710 920-1030"
    181 = "2930-3085"
    186 = "285-370"
    187 = "2420-2680 2830-2880 3210-3300 4140-4290"
    188 = "1235-1400 1515-1580"
    189 = "2540-2640 2860-2945 3190-3270"
    194 = "-240"
}

foreach ($row in $hValues.Keys) {
    $cell = $ws.Range("H$row")
    $cell.NumberFormat = "@"
    $cell.Value = $hValues[$row]
}

# --- Formatting touch-up on E155 -------------------------------------------
# The label cell gets switched onto an explicit (non-default) font run,
# introducing a new font in the workbook's style table.
$ws.Range("E155").Font.Name = "Calibri"

# --- Leave the view where the user ended up editing ------------------------
$ws.Range("D192").Select()
